$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (model_name, precision, recall, f1_score,
#           %_correctly_detected_entities, %_correctly_detected_and_classified_entities,
#           %_exact_matches_start_pos, %_within_2_matches_start_pos, %_within_5_matches_start_pos, %_within_10_matches_start_pos,
#           %_exact_matches_end_pos, %_within_2_matches_end_pos, %_within_5_matches_end_pos, %_within_10_matches_end_pos)
# for rows 2-9, reflecting the new multiprocessing run outputs.

$data = @(
    @("fine_tuned_gpt_4o_mini",    0.6507966344243519, 0.6455199314574315, 0.6364691594862174, 65.21739130434783, 56.52173913043478, 16.90476190476191, 28.09523809523809, 36.19047619047619, 45.95238095238095, 15.47619047619048, 27.61904761904762, 35.71428571428572, 45.95238095238095),
    @("fine_tuned_gpt_4o",         0.6843809432830815, 0.6985712782587784, 0.6808901394872756, 72.67080745341616, 64.13043478260869, 19.87179487179487, 35.04273504273504, 50.85470085470085, 62.39316239316239, 17.73504273504274, 34.4017094017094,  50.85470085470085, 62.17948717948718),
    @("fine_tuned_gpt_4o_mini_v2", 0.8562024897414151, 0.8599851537351538, 0.8538762209899445, 84.93788819875776, 75.77639751552795, 21.57221206581353, 34.55210237659963, 49.90859232175503, 61.42595978062158, 20.84095063985375, 34.36928702010969, 49.72577696526508, 61.24314442413162),
    @("fine_tuned_gpt_4o_v2",      0.8102967337547931, 0.8163485675204426, 0.8040474712754717, 81.67701863354037, 75.62111801242236, 30.03802281368822, 45.81749049429658, 60.45627376425855, 69.01140684410646, 30.03802281368822, 45.81749049429658, 60.26615969581749, 69.01140684410646),
    @("gpt_4o",                    0.7619155036301417, 0.8016875919219669, 0.7684494353043774, 74.04921700223713, 45.63758389261745, 25.98187311178248, 39.57703927492447, 51.05740181268882, 59.21450151057401, 21.45015105740181, 37.76435045317221, 51.3595166163142,  59.21450151057401),
    @("gliner",                    0.6532512626262627, 0.4476419992044993, 0.506185989190214,  51.24223602484472, 20.3416149068323,  98.7878787878788,  98.7878787878788,  98.7878787878788,  99.09090909090909, 98.7878787878788,  98.7878787878788,  98.7878787878788,  99.09090909090909),
    @("gpt_4o_mini",               0.5237153074477313, 0.6889567174723424, 0.5736679606668337, 68.63354037267081, 39.59627329192546, 11.08597285067873, 18.09954751131222, 22.85067873303167, 28.28054298642534, 6.108597285067873, 16.7420814479638,  23.30316742081448, 27.14932126696833),
    @("sonnet_35",                 0.6317338963025781, 0.7606979392135642, 0.6731942184922062, 76.86335403726709, 67.2360248447205,  23.83838383838384, 38.98989898989899, 51.71717171717172, 62.22222222222222, 20,                 38.98989898989899, 51.71717171717172, 62.02020202020202)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
